$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.191131666666666
$ws.Range("H2").Value = 3.573395
$ws.Range("I2").Value = 0.02720036629735778
$ws.Range("J2").Value = 0.02720036629735778
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 117.044563
$ws.Range("N2").Value = 351.133689
$ws.Range("O2").Value = 0.3245365645427815
$ws.Range("P2").Value = 0.3245365645427815
$ws.Range("Q2").Value = 139.4154854004616
$ws.Range("R2").Value = 1254.739368604155
$ws.Range("S2").Value = 0.00882751343244975
$ws.Range("T2").Value = 0.00882751343244975
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.191131666666666
$ws.Range("H3").Value = 3.573395
$ws.Range("I3").Value = 0.02720036629735778
$ws.Range("J3").Value = 0.02720036629735778
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 101.5800373333333
$ws.Range("N3").Value = 304.740112
$ws.Range("O3").Value = 0.281657135515876
$ws.Range("P3").Value = 0.281657135515876
$ws.Range("Q3").Value = 120.9951991689155
$ws.Range("R3").Value = 1088.95679252024
$ws.Range("S3").Value = 0.007661177256296366
$ws.Range("T3").Value = 0.007661177256296364
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.191131666666666
$ws.Range("H4").Value = 3.573395
$ws.Range("I4").Value = 0.02720036629735778
$ws.Range("J4").Value = 0.02720036629735778
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 142.0267893333333
$ws.Range("N4").Value = 426.080368
$ws.Range("O4").Value = 0.3938062999413425
$ws.Range("P4").Value = 0.3938062999413425
$ws.Range("Q4").Value = 169.1726062899289
$ws.Range("R4").Value = 1522.55345660936
$ws.Range("S4").Value = 0.01071167560861166
$ws.Range("T4").Value = 0.01071167560861166
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 34.415161
$ws.Range("H5").Value = 103.245483
$ws.Range("I5").Value = 0.7858954736735307
$ws.Range("J5").Value = 0.7858954736735306
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 117.044563
$ws.Range("N5").Value = 351.133689
$ws.Range("O5").Value = 0.3245365645427815
$ws.Range("P5").Value = 0.3245365645427815
$ws.Range("Q5").Value = 4028.107479819643
$ws.Range("R5").Value = 36252.96731837679
$ws.Range("S5").Value = 0.2550518171157296
$ws.Range("T5").Value = 0.2550518171157296
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 34.415161
$ws.Range("H6").Value = 103.245483
$ws.Range("I6").Value = 0.7858954736735307
$ws.Range("J6").Value = 0.7858954736735306
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 101.5800373333333
$ws.Range("N6").Value = 304.740112
$ws.Range("O6").Value = 0.281657135515876
$ws.Range("P6").Value = 0.281657135515876
$ws.Range("Q6").Value = 3495.893339212678
$ws.Range("R6").Value = 31463.0400529141
$ws.Range("S6").Value = 0.2213530679297792
$ws.Range("T6").Value = 0.2213530679297791
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 34.415161
$ws.Range("H7").Value = 103.245483
$ws.Range("I7").Value = 0.7858954736735307
$ws.Range("J7").Value = 0.7858954736735306
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 142.0267893333333
$ws.Range("N7").Value = 426.080368
$ws.Range("O7").Value = 0.3938062999413425
$ws.Range("P7").Value = 0.3938062999413425
$ws.Range("Q7").Value = 4887.874821219751
$ws.Range("R7").Value = 43990.87339097775
$ws.Range("S7").Value = 0.3094905886280219
$ws.Range("T7").Value = 0.3094905886280219
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.184723
$ws.Range("H8").Value = 24.554169
$ws.Range("I8").Value = 0.1869041600291116
$ws.Range("J8").Value = 0.1869041600291116
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 117.044563
$ws.Range("N8").Value = 351.133689
$ws.Range("O8").Value = 0.3245365645427815
$ws.Range("P8").Value = 0.3245365645427815
$ws.Range("Q8").Value = 957.977326811049
$ws.Range("R8").Value = 8621.795941299442
$ws.Range("S8").Value = 0.06065723399460213
$ws.Range("T8").Value = 0.06065723399460213
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.184723
$ws.Range("H9").Value = 24.554169
$ws.Range("I9").Value = 0.1869041600291116
$ws.Range("J9").Value = 0.1869041600291116
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 101.5800373333333
$ws.Range("N9").Value = 304.740112
$ws.Range("O9").Value = 0.281657135515876
$ws.Range("P9").Value = 0.281657135515876
$ws.Range("Q9").Value = 831.404467902992
$ws.Range("R9").Value = 7482.640211126929
$ws.Range("S9").Value = 0.05264289032980047
$ws.Range("T9").Value = 0.05264289032980045
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.184723
$ws.Range("H10").Value = 24.554169
$ws.Range("I10").Value = 0.1869041600291116
$ws.Range("J10").Value = 0.1869041600291116
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 142.0267893333333
$ws.Range("N10").Value = 426.080368
$ws.Range("O10").Value = 0.3938062999413425
$ws.Range("P10").Value = 0.3938062999413425
$ws.Range("Q10").Value = 1162.449929272688
$ws.Range("R10").Value = 10462.04936345419
$ws.Range("S10").Value = 0.07360403570470901
$ws.Range("T10").Value = 0.073604035704709
